$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the DenseNet122 / DenseNet123 model rows to DenseNet121 (merging them into the
# existing DenseNet121 label). This also causes Excel to drop the now-unused
# "DenseNet122"/"DenseNet123" shared-string entries on save.
$ws.Range("A36").Value = "DenseNet121"
$ws.Range("A37").Value = "DenseNet121"

# Hide the rows that are no longer the focus of the pipeline run (rows 2-25)
$ws.Range("A2:A25").EntireRow.Hidden = $true

# Update the recomputed RNN / STL10 attack-model metrics
$ws.Range("K26").Value = 0.64383012056350697
$ws.Range("L26").Value = 0.58931088813334398
$ws.Range("M26").Value = 0.949051816239316
$ws.Range("N26").Value = 0.33860844017093999
$ws.Range("O26").Value = 0.72711925103596398

$ws.Range("K27").Value = 0.64072513580322199
$ws.Range("L27").Value = 0.587970113119338
$ws.Range("M27").Value = 0.94057158119658102
$ws.Range("N27").Value = 0.34087873931623902
$ws.Range("O27").Value = 0.72360208563429396

$ws.Range("K30").Value = 0.65064102411270097
$ws.Range("L30").Value = 0.59598366235534295
$ws.Range("M30").Value = 0.93536324786324698
$ws.Range("N30").Value = 0.36591880341880301
$ws.Range("O30").Value = 0.72806652806652805

$ws.Range("K31").Value = 0.64232772588729803
$ws.Range("L31").Value = 0.595024742543801
$ws.Range("M31").Value = 0.89122596153846101
$ws.Range("N31").Value = 0.393429487179487
$ws.Range("O31").Value = 0.713609752185419

$ws.Range("K34").Value = 0.60803955793380704
$ws.Range("L34").Value = 0.57156758669497498
$ws.Range("M34").Value = 0.86284722222222199
$ws.Range("N34").Value = 0.35323183760683702
$ws.Range("O34").Value = 0.68763303533418396

$ws.Range("K35").Value = 0.62937366962432795
$ws.Range("L35").Value = 0.57617456261057598
$ws.Range("M35").Value = 0.97856570512820495
$ws.Range("N35").Value = 0.28018162393162299
$ws.Range("O35").Value = 0.72529756749399898

$ws.Range("K38").Value = 0.728565692901611
$ws.Range("L38").Value = 0.65077966698969203
$ws.Range("M38").Value = 0.98651175213675202
$ws.Range("N38").Value = 0.470619658119658
$ws.Range("O38").Value = 0.78422421572270296

$ws.Range("K39").Value = 0.73150372505187899
$ws.Range("L39").Value = 0.65312251567882695
$ws.Range("M39").Value = 0.98744658119658102
$ws.Range("N39").Value = 0.47556089743589702
$ws.Range("O39").Value = 0.78621936307087004

$ws.Range("K42").Value = 0.70382612943649203
$ws.Range("L42").Value = 0.63283867879368105
$ws.Range("M42").Value = 0.97102029914529897
$ws.Range("N42").Value = 0.43663194444444398
$ws.Range("O42").Value = 0.76627585298379597

$ws.Range("K43").Value = 0.71624600887298495
$ws.Range("L43").Value = 0.64026159643120095
$ws.Range("M43").Value = 0.98711271367521303
$ws.Range("N43").Value = 0.44537927350427298
$ws.Range("O43").Value = 0.776724024694601

$ws.Range("K46").Value = 0.68676549196243197
$ws.Range("L46").Value = 0.64402677651905205
$ws.Range("M46").Value = 0.83513621794871795
$ws.Range("N46").Value = 0.53839476495726402
$ws.Range("O46").Value = 0.727235725084312

$ws.Range("K47").Value = 0.69140625
$ws.Range("L47").Value = 0.64691712367382503
$ws.Range("M47").Value = 0.84281517094017
$ws.Range("N47").Value = 0.539997329059829
$ws.Range("O47").Value = 0.73198596572621499

# Update the view: scroll/selection moved to I51 (and the window no longer needs to
# stay pinned at A35 now that rows 2-25 are hidden).
$ws.Range("I51").Select()
